$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G2 / H2: stable coin symbol + address swapped from USDC to USDT (BNB chain) ---
$ws.Range("G2").Value = "USDT"
$ws.Range("H2").Value = "0x55d398326f99059ff775485246999027b3197955"

# --- C2: target gain percent lowered from 20 to 10 ---
$ws.Range("C2").Value = 10

# --- Column B width tightened ---
$ws.Columns.Item(2).ColumnWidth = 17.8

# --- Row 4 (existing first data row) gets new price / qty / date ---
$ws.Range("A4").Value = 0.2675
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 45921

# --- Row 5 ---
$ws.Range("A5").Value = 0.2427
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 45921
$ws.Range("C5").NumberFormat = "dd/mm/yy"

# --- Row 6 ---
$ws.Range("A6").Value = 0.2282
$ws.Range("B6").Value = 10

# --- Row 7 ---
$ws.Range("A7").Value = 0.18483
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "31/10/2025"
$ws.Range("D7").Value = "18:45:59"

# --- Row 8 (new) ---
$ws.Range("A8").Value = 0.18309
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "02/11/2025"
$ws.Range("D8").Value = "21:02:51"

# --- Row 9 (new) ---
$ws.Range("A9").Value = 0.16994
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "03/11/2025"
$ws.Range("D9").Value = "19:32:04"

# --- Selection moved to B1 ---
$ws.Range("B1").Select()
